# Update COVID-19 countries data and Spain provincias timestamp
# Generated to reflect: reorder Irlanda/Suecia/India, reorder Jordania/Malta/Taiwan/Reunion,
# refresh several countries' case counts, and bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 14 de Abril de 2020 a las 19:22'
$ws.Cells.Item(4, 2).Value = 598737
$ws.Cells.Item(4, 3).Value = 11796
$ws.Cells.Item(4, 4).Value = 38015
$ws.Cells.Item(4, 5).Value = 535952
$ws.Cells.Item(4, 7).Value = 1130
$ws.Cells.Item(4, 8).Value = 24770
$ws.Cells.Item(12, 2).Value = 65111
$ws.Cells.Item(12, 3).Value = 4062
$ws.Cells.Item(12, 4).Value = 4799
$ws.Cells.Item(12, 5).Value = 58909
$ws.Cells.Item(12, 6).Value = 1809
$ws.Cells.Item(12, 7).Value = 107
$ws.Cells.Item(12, 8).Value = 1403
$ws.Cells.Item(15, 4).Value = 8008
$ws.Cells.Item(15, 5).Value = 17375
$ws.Cells.Item(22, 1).Value = 'Irlanda'
$ws.Cells.Item(22, 2).Value = 11479
$ws.Cells.Item(22, 3).Value = 832
$ws.Cells.Item(22, 4).Value = 25
$ws.Cells.Item(22, 5).Value = 11048
$ws.Cells.Item(22, 6).Value = 194
$ws.Cells.Item(22, 7).Value = 41
$ws.Cells.Item(22, 8).Value = 406
$ws.Cells.Item(23, 1).Value = 'Suecia'
$ws.Cells.Item(23, 2).Value = 11445
$ws.Cells.Item(23, 3).Value = 497
$ws.Cells.Item(23, 4).Value = 381
$ws.Cells.Item(23, 5).Value = 10031
$ws.Cells.Item(23, 6).Value = 915
$ws.Cells.Item(23, 7).Value = 114
$ws.Cells.Item(23, 8).Value = 1033
$ws.Cells.Item(24, 1).Value = 'India'
$ws.Cells.Item(24, 2).Value = 10941
$ws.Cells.Item(24, 3).Value = 488
$ws.Cells.Item(24, 4).Value = 1295
$ws.Cells.Item(24, 5).Value = 9278
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 10
$ws.Cells.Item(24, 8).Value = 368
$ws.Cells.Item(29, 2).Value = 7603
$ws.Cells.Item(29, 3).Value = 74
$ws.Cells.Item(29, 5).Value = 6651
$ws.Cells.Item(35, 2).Value = 6111
$ws.Cells.Item(35, 3).Value = 52
$ws.Cells.Item(35, 5).Value = 5308
$ws.Cells.Item(35, 6).Value = 428
$ws.Cells.Item(60, 2).Value = 1888
$ws.Cells.Item(60, 3).Value = 125
$ws.Cells.Item(60, 4).Value = 217
$ws.Cells.Item(60, 5).Value = 1545
$ws.Cells.Item(100, 1).Value = 'Jordania'
$ws.Cells.Item(100, 2).Value = 397
$ws.Cells.Item(100, 3).Value = 6
$ws.Cells.Item(100, 4).Value = 235
$ws.Cells.Item(100, 5).Value = 155
$ws.Cells.Item(100, 6).Value = 5
$ws.Cells.Item(100, 8).Value = 7
$ws.Cells.Item(101, 1).Value = 'Malta'
$ws.Cells.Item(101, 3).Value = 9
$ws.Cells.Item(101, 4).Value = 44
$ws.Cells.Item(101, 5).Value = 346
$ws.Cells.Item(101, 6).Value = 4
$ws.Cells.Item(101, 8).Value = 3
$ws.Cells.Item(102, 1).Value = 'Taiwan'
$ws.Cells.Item(102, 2).Value = 393
$ws.Cells.Item(102, 4).Value = 124
$ws.Cells.Item(102, 5).Value = 263
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 8).Value = 6
$ws.Cells.Item(103, 1).Value = 'Reunion'
$ws.Cells.Item(103, 2).Value = 393
$ws.Cells.Item(103, 4).Value = 40
$ws.Cells.Item(103, 5).Value = 351
$ws.Cells.Item(103, 6).Value = 3
$ws.Cells.Item(103, 8).Value = 0
